$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 37499.5
$ws.Range("J3").Value = 37499.5
$ws.Range("L3").Value = 37499.5
$ws.Range("N3").Value = -37727.5
$ws.Range("H5").Value = 189.58824
$ws.Range("I5").Value = 97.15385000000001
$ws.Range("J5").Value = 490
$ws.Range("K5").Value = 97.15385000000001
$ws.Range("L5").Value = 490
$ws.Range("M5").Value = 17.84614999999999
$ws.Range("N5").Value = -720
$ws.Range("H17").Value = 1405515.8
$ws.Range("J17").Value = 1778052.6
$ws.Range("L17").Value = 5334157.800000001
$ws.Range("N17").Value = -5334493.800000001
$ws.Range("H102").Value = 37499.5
$ws.Range("J102").Value = 37499.5
$ws.Range("L102").Value = 37499.5
$ws.Range("N102").Value = -43989.5
$ws.Range("H113").Value = 7997.5
$ws.Range("I113").Value = 7995
$ws.Range("K113").Value = 7995
$ws.Range("M113").Value = -4741
$ws.Range("H129").Value = 3239.2273
$ws.Range("J129").Value = 6459.4
$ws.Range("L129").Value = 19378.2
$ws.Range("N129").Value = -29378.2
$ws.Range("H132").Value = 1247
$ws.Range("I132").Value = 1306.0526
$ws.Range("K132").Value = 3918.1578
$ws.Range("M132").Value = -1388.1578
$ws.Range("H135").Value = 858
$ws.Range("J135").Value = 274.5
$ws.Range("L135").Value = 2470.5
$ws.Range("N135").Value = -7540.5
$ws.Range("H137").Value = 9733.816000000001
$ws.Range("I137").Value = 4985.147
$ws.Range("J137").Value = 15943.615
$ws.Range("K137").Value = 14955.441
$ws.Range("L137").Value = 47830.845
$ws.Range("M137").Value = -12405.441
$ws.Range("N137").Value = -52930.845
$ws.Range("H138").Value = 2442.95
$ws.Range("I138").Value = 796.7692
$ws.Range("K138").Value = 2390.3076
$ws.Range("M138").Value = 2749.6924

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3534.1516
$ws.Range("I32").Value = 2758.3442
$ws.Range("K32").Value = 2758.3442
$ws.Range("M32").Value = -2471.3442
$ws.Range("H45").Value = 11334.643
$ws.Range("I45").Value = 11334.643
$ws.Range("K45").Value = 11334.643
$ws.Range("M45").Value = -10957.643
$ws.Range("H61").Value = 11511.267
$ws.Range("I61").Value = 4538.222
$ws.Range("K61").Value = 4538.222
$ws.Range("M61").Value = -4326.222
$ws.Range("H74").Value = 15335.044
$ws.Range("I74").Value = 20655.066
$ws.Range("K74").Value = 20655.066
$ws.Range("M74").Value = -19781.066
$ws.Range("H77").Value = 15335.044
$ws.Range("I77").Value = 20655.066
$ws.Range("K77").Value = 103275.33
$ws.Range("M77").Value = -98907.32999999999
$ws.Range("H92").Value = 49999
$ws.Range("J92").Value = 49999
$ws.Range("L92").Value = 49999
$ws.Range("N92").Value = -54991
$ws.Range("H93").Value = 39250
$ws.Range("I93").Value = 39000
$ws.Range("J93").Value = 39500
$ws.Range("K93").Value = 39000
$ws.Range("L93").Value = 39500
$ws.Range("M93").Value = -36504
$ws.Range("N93").Value = -44492
$ws.Range("H122").Value = 3707.1428
$ws.Range("I122").Value = 4950
$ws.Range("J122").Value = 3500
$ws.Range("K122").Value = 14850
$ws.Range("L122").Value = 10500
$ws.Range("M122").Value = -12400
$ws.Range("N122").Value = -15400
$ws.Range("H132").Value = 6682.909
$ws.Range("I132").Value = 5520.5625
$ws.Range("J132").Value = 9782.5
$ws.Range("K132").Value = 16561.6875
$ws.Range("L132").Value = 29347.5
$ws.Range("M132").Value = -14031.6875
$ws.Range("N132").Value = -34407.5
$ws.Range("H136").Value = 11511.267
$ws.Range("I136").Value = 4538.222
$ws.Range("K136").Value = 13614.666
$ws.Range("M136").Value = -11064.666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H16").Value = 6877.1665
$ws.Range("I16").Value = 6862.6
$ws.Range("J16").Value = 6950
$ws.Range("K16").Value = 6862.6
$ws.Range("L16").Value = 6950
$ws.Range("M16").Value = -6692.6
$ws.Range("N16").Value = -7290
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H95").Value = 9312
$ws.Range("J95").Value = 9312
$ws.Range("L95").Value = 9312
$ws.Range("N95").Value = -14804
$ws.Range("H97").Value = 5500
$ws.Range("I97").Value = 5500
$ws.Range("K97").Value = 5500
$ws.Range("M97").Value = -4509
$ws.Range("H100").Value = 45484.125
$ws.Range("J100").Value = 45484.125
$ws.Range("L100").Value = 45484.125
$ws.Range("N100").Value = -47648.125
$ws.Range("H134").Value = 10173.917
$ws.Range("I134").Value = 5691
$ws.Range("J134").Value = 16450
$ws.Range("K134").Value = 17073
$ws.Range("L134").Value = 49350
$ws.Range("M134").Value = -14538
$ws.Range("N134").Value = -54420

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 88297144
$ws.Range("I4").Value = 22498.584
$ws.Range("K4").Value = 22498.584
$ws.Range("M4").Value = -22386.584
$ws.Range("H7").Value = 353.63333
$ws.Range("I7").Value = 377.58334
$ws.Range("K7").Value = 377.58334
$ws.Range("M7").Value = -264.58334
$ws.Range("H31").Value = 358278.56
$ws.Range("I31").Value = 73154.86
$ws.Range("J31").Value = 593086.3
$ws.Range("K31").Value = 73154.86
$ws.Range("L31").Value = 593086.3
$ws.Range("M31").Value = -72859.86
$ws.Range("N31").Value = -593676.3
$ws.Range("H34").Value = 358278.56
$ws.Range("I34").Value = 73154.86
$ws.Range("J34").Value = 593086.3
$ws.Range("K34").Value = 73154.86
$ws.Range("L34").Value = 593086.3
$ws.Range("M34").Value = -72952.86
$ws.Range("N34").Value = -593490.3
$ws.Range("H132").Value = 31536.693
$ws.Range("J132").Value = 49902.707
$ws.Range("L132").Value = 149708.121
$ws.Range("N132").Value = -154768.121
$ws.Range("H134").Value = 1841.3793
$ws.Range("I134").Value = 1371.4642
$ws.Range("K134").Value = 4114.392599999999
$ws.Range("M134").Value = -1579.392599999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 1999.5
$ws.Range("I3").Value = 1999.5
$ws.Range("K3").Value = 5998.5
$ws.Range("M3").Value = -5886.5
$ws.Range("H114").Value = 389
$ws.Range("J114").Value = 378
$ws.Range("L114").Value = 1134
$ws.Range("N114").Value = -7642
$ws.Range("H129").Value = 2956.5833
$ws.Range("I129").Value = 1207.6666
$ws.Range("K129").Value = 3622.9998
$ws.Range("M129").Value = 1377.0002
$ws.Range("H139").Value = 1000
$ws.Range("I139").Value = 1000
$ws.Range("K139").Value = 3000
$ws.Range("M139").Value = 2140

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 11499.5
$ws.Range("I5").Value = 8999.666999999999
$ws.Range("K5").Value = 8999.666999999999
$ws.Range("M5").Value = -8887.666999999999
$ws.Range("H42").Value = 65712
$ws.Range("I42").Value = 59962
$ws.Range("J42").Value = 67628.664
$ws.Range("K42").Value = 59962
$ws.Range("L42").Value = 67628.664
$ws.Range("M42").Value = -59477
$ws.Range("N42").Value = -68598.664
$ws.Range("H115").Value = 65712
$ws.Range("I115").Value = 59962
$ws.Range("J115").Value = 67628.664
$ws.Range("K115").Value = 59962
$ws.Range("L115").Value = 67628.664
$ws.Range("M115").Value = -58787
$ws.Range("N115").Value = -69978.664
$ws.Range("H122").Value = 7734.9116
$ws.Range("I122").Value = 5307.5
$ws.Range("J122").Value = 9058.954
$ws.Range("K122").Value = 15922.5
$ws.Range("L122").Value = 27176.862
$ws.Range("M122").Value = -13472.5
$ws.Range("N122").Value = -32076.862
$ws.Range("H132").Value = 24876.215
$ws.Range("I132").Value = 24356.8
$ws.Range("K132").Value = 73070.39999999999
$ws.Range("M132").Value = -70540.39999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 735.3333
$ws.Range("I22").Value = 544.6
$ws.Range("K22").Value = 544.6
$ws.Range("M22").Value = -249.6
$ws.Range("H27").Value = 735.3333
$ws.Range("I27").Value = 544.6
$ws.Range("K27").Value = 544.6
$ws.Range("M27").Value = -437.6
$ws.Range("H132").Value = 3510.8071
$ws.Range("I132").Value = 2920.359
$ws.Range("J132").Value = 4790.1113
$ws.Range("K132").Value = 8761.076999999999
$ws.Range("L132").Value = 14370.3339
$ws.Range("M132").Value = -6231.076999999999
$ws.Range("N132").Value = -19430.3339
$ws.Range("H136").Value = 2633.5557
$ws.Range("I136").Value = 2286.85
$ws.Range("K136").Value = 6860.549999999999
$ws.Range("M136").Value = -4310.549999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 3000.5
$ws.Range("I2").Value = 3000.5
$ws.Range("K2").Value = 3000.5
$ws.Range("M2").Value = -2888.5
$ws.Range("H136").Value = 2314.1702
$ws.Range("I136").Value = 768.7308
$ws.Range("J136").Value = 4227.5713
$ws.Range("K136").Value = 2306.1924
$ws.Range("L136").Value = 12682.7139
$ws.Range("M136").Value = 243.8076000000001
$ws.Range("N136").Value = -17782.7139
